$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the "FAR Case 2018-004" and ")" runs (identical formatting) into a
#    single run reading "FAR Case 2018-004)". A Find/Replace over the exact
#    same text naturally coalesces the two adjacent, identically-formatted
#    runs into one run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("FAR Case 2018-004)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "FAR Case 2018-004)", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Split "Effective August 31, 2020" into three runs: "Effective", ":" and
#    " August 31, 2020". Plain text insertion/replacement gets coalesced by
#    the engine into a single run whenever the resulting runs would have
#    identical formatting, so instead we replace the whole paragraph (via
#    InsertXML on its Range) with explicit separate <w:r> elements.
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$p2Start = $p2.Range.Start
$p2End = $p2.Range.End
$p2Range = $d.Range($p2Start, $p2End)
$effectiveXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:pPr><w:pStyle w:val="Heading2"/></w:pPr>' + `
    '<w:r><w:t>Effective</w:t></w:r>' + `
    '<w:r><w:t>:</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> August 31, 2020</w:t></w:r>' + `
    '</w:p>'
$p2Range.InsertXML($effectiveXml)

# ---------------------------------------------------------------------------
# 3) Change the paragraph-mark formatting of the paragraph that contains the
#    first hyperlink (rId4) from rStyle="Hyperlink" to explicit direct
#    formatting (blue color + single underline), without disturbing the
#    hyperlink run itself.
#
#    InsertXML silently keeps a paragraph's own pPr when the inserted
#    fragment is absorbed as a single paragraph, so instead we insert two
#    empty paragraphs right after the existing one (second one carrying the
#    desired pPr) and then delete the two intervening paragraph marks so
#    that the desired pPr "flows back" onto the original paragraph while its
#    run content is left completely untouched.
# ---------------------------------------------------------------------------
$pLink = $d.Paragraphs.Item(4)
$linkEnd = $pLink.Range.End
$insertPoint = $d.Range($linkEnd, $linkEnd)
$markXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>' + `
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:pPr><w:rPr><w:color w:val="0000FF"/><w:u w:val="single"/></w:rPr></w:pPr>' + `
    '</w:p>'
$insertPoint.InsertXML($markXml)

$pA = $d.Paragraphs.Item(4)
$pB = $d.Paragraphs.Item(5)
$d.Range($pA.Range.End - 1, $pB.Range.End - 1).Delete() | Out-Null

$pA2 = $d.Paragraphs.Item(4)
$pB2 = $d.Paragraphs.Item(5)
$d.Range($pA2.Range.End - 1, $pB2.Range.End - 1).Delete() | Out-Null

# ---------------------------------------------------------------------------
# 4) Remove the trailing paragraphs: the empty Hyperlink-styled paragraph,
#    the "Start here: <link>" paragraph, and the final empty paragraph. They
#    are all the remaining content after the (now-reformatted) hyperlink
#    paragraph, right up to the end of the document body's text.
# ---------------------------------------------------------------------------
$pLinkFinal = $d.Paragraphs.Item(4)
$tailRange = $d.Range($pLinkFinal.Range.End, $d.Content.End)
$tailRange.Delete() | Out-Null

Write-Output "done"
